# Lägg till Vindskiveplåt - beräknas i antal (st) baserat på täckning per plåt
#
# Insert a new row above the existing "Takfotsbräda 22x145" row (row 20)
# on the "Material" sheet, shifting it and the following rows down by one,
# then fill the new row with the "Vindskiveplåt" article data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Material")

# Insert a new row at row 20; existing rows 20-22 shift down to 21-23.
$ws.Rows.Item(20).Insert()

# Populate the newly inserted row 20 with the Vindskiveplåt data.
$ws.Range("A20").Value = "Tak"
$ws.Range("B20").Value = "Vindskiveplåt"
$ws.Range("C20").Value = "st"
$ws.Range("D20").Value = 1.9
$ws.Range("E20").Value = 0.15
$ws.Range("F20").Value = 85
$ws.Range("G20").Value = 145
$ws.Range("H20").Value = $true
$ws.Range("I20").Value = "Täcker 1.9m/st, avrundas uppåt"
